# Applies the "black holes" -> "biology" content rewrite described by the
# commit diff. All substitutions locate the target text with Find (without
# using Find's own Replace argument, which runs the text through
# AutoCorrect/AutoFormat and mangles straight quotes into curly ones) and
# then set Range.Text directly, which preserves the run's existing
# formatting (rFonts/color/sz) and leaves apostrophes untouched. New
# sentences are appended with InsertAfter using the same rPr as their
# neighbouring run.

$d = $word.ActiveDocument

function Find-Range($old) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $old"
    }
    return $rng
}

function Replace-Text($old, $new) {
    $rng = Find-Range $old
    $rng.Text = $new
}

# --- Title -------------------------------------------------------------
Replace-Text "Unraveling the Enigma of Black Holes" "Biology: Exploring the Symphony of Life"

# --- Byline --------------------------------------------------------------
Replace-Text " Amelia Carter" " Emily Jones"

# --- Email: "ameliacarter@cosmosinstitute" -> "emily" + "." + "jones@eduworld" ---
# (the "." and "org" that already follow the replaced run become
#  "jones@eduworld" + existing "." + existing "org" -> "jones@eduworld.org")
Replace-Text "ameliacarter@cosmosinstitute" "emily.jones@eduworld"

# --- Paragraph 1 (intro) ------------------------------------------------
Replace-Text "In the boundless expanse of the universe, there lies a celestial enigma that has captivated the minds of scientists and lay people alike: black holes" "Biology, the study of life, is a journey that unravels the mysteries of the living world"

Replace-Text " These enigmatic entities are cosmic vacuums with an infinitely strong gravitational pull, from which nothing, not even light, can escape" " It's an exploration into the intricate workings of organisms, from the smallest microbes to the largest whales, and the dynamic interactions between them"

Replace-Text " They represent a dark frontier in our understanding of space and time, inviting us to explore their extraordinary properties and unravel their cosmic mysteries" " We embark on a quest to comprehend the secrets of life, delving into the symphony of biological processes that orchestrate the harmony of existence"

# New sentences inserted after "...harmony of existence"
$anchor = Find-Range "harmony of existence"
$ins = $d.Range($anchor.End, $anchor.End)
$ins.InsertAfter(".")
$ins.Collapse(0)
$ins.InsertAfter(" As we dissect the molecular mechanisms that govern cellular functions, we unravel the mysteries of genetics, the blueprint of life")
$ins.Collapse(0)
$ins.InsertAfter(".")
$ins.Collapse(0)
$ins.InsertAfter(" We investigate the intricate web of ecosystems, the intricate relationships between organisms and their environments, and the delicate balance that sustains our planet's biodiversity")

# --- Paragraph 1 continued (after the <br/><br/>) -----------------------
Replace-Text "In this exploration, we embark on a journey to penetrate the veil of darkness surrounding black holes" "Biology propels us into the realm of human biology, unraveling the complexities of our bodies, from the microscopic world of cells to the sophisticated systems that regulate our thoughts, emotions, and actions"

Replace-Text " We begin by unraveling their gravitational peculiarities, delving into the concept of event horizons and their role in creating an inescapable boundary" " We delve into the wonders of the immune system, our body's intricate defense mechanism, fending off invaders and maintaining internal harmony"

Replace-Text " Additionally, we examine the mind-bending phenomena occurring near black holes, including time dilation and the mesmerizing behavior of light" " We explore the marvels of reproduction, the miracle of new life emerging from the fusion of genetic material"

# New sentences inserted after "...fusion of genetic material"
$anchor2 = Find-Range "fusion of genetic material"
$ins2 = $d.Range($anchor2.End, $anchor2.End)
$ins2.InsertAfter(".")
$ins2.Collapse(0)
$ins2.InsertAfter(" By studying the human body, we gain a profound appreciation for the resilience, adaptability, and interconnectedness of life")

# --- Paragraph 1 continued (after the second <br/><br/>) -----------------
Replace-Text "Furthermore, we investigate the birth of these cosmic leviathans through the death of massive stars" "Furthermore, biology illuminates the interconnectedness of life on Earth"

Replace-Text " We analyze the various evolutionary pathways, such as stellar collapse and supernova explosions, that lead to the formation of these enigmatic entities" " We uncover the intricate web of interdependence among organisms, the delicate balance of ecosystems, and the profound impact of human activities on the natural world"

Replace-Text " Our quest for knowledge leads us to question the ultimate fate of black holes, considering their hypothetical evaporation through Hawking radiation and their possible involvement in mysterious cosmic phenomena like gravitational waves" " Biology empowers us to recognize our responsibility as stewards of the planet, inspiring us to act as conscientious citizens, preserving and protecting the diversity of life for generations to come"

# --- Summary paragraph ----------------------------------------------------
Replace-Text "Black holes stand as a testament to the vastness and complexity of the universe, beckoning us to push the boundaries of our scientific understanding" "Biology is a captivating subject that delves into the intricacies of life, unraveling the mysteries of living organisms and their interactions with each other and their environment"

Replace-Text " Through an examination of their gravitational anomalies, formation mechanisms, and cosmic interactions, we have delved into the depths of these celestial conundrums" " By exploring the symphony of biological processes, we gain a deeper understanding of our bodies, appreciate the marvels of life's diversity, and recognize our role as stewards of the planet"

Replace-Text " As we continue to unravel the enigma of black holes, we unlock new insights into the fundamental nature of space, time, and the dynamics of the cosmos" " Biology inspires us to question, investigate, and marvel at the wonders of the natural world, nurturing a lifelong appreciation for the beauty and complexity of life"

# --- Trailing empty paragraph ---------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Output $d.Content.Text
